# "Generate Report for Archive"
# The status "Ready for handoff" has moved on to "In Translation" for both
# source documents, across the Overview sheet (per-language status columns)
# and each per-language detail sheet (zh-cn / de-de "Status" column). Excel
# auto-shrinks those status columns afterwards, since the new text is
# shorter than the old text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: per-language status columns E (zh-cn) and F (de-de)
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Detail sheets: column C is "Status"
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# Column widths shrink to fit the shorter status text
$overview.Columns("E:F").ColumnWidth = 12.5
$zhcn.Columns("C:C").ColumnWidth = 12.5
$dede.Columns("C:C").ColumnWidth = 12.5
